$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1514.7329694126095
$ws.Range("E2").Value = 1731.0115884581783
$ws.Range("G2").Value = 957.2368372616152
$ws.Range("J2").Value = 471.2227821772243
$ws.Range("C3").Value = 1060.1522125646386
$ws.Range("E3").Value = 1813.0016146387984
$ws.Range("G3").Value = 1137.1588884089708
$ws.Range("J3").Value = 537.1585582369928
$ws.Range("C4").Value = 1089.1068138559324
$ws.Range("E4").Value = 1791.6388159567416
$ws.Range("G4").Value = 825.6203014835611
$ws.Range("J4").Value = 880.0915727635413
$ws.Range("C5").Value = 1085.5815223093905
$ws.Range("E5").Value = 1821.7586330497975
$ws.Range("G5").Value = 1082.9627095654841
$ws.Range("J5").Value = 713.237348637963
$ws.Range("C6").Value = 659.9592090102584
$ws.Range("E6").Value = 1908.8819777038418
$ws.Range("G6").Value = 805.8749150451886
$ws.Range("J6").Value = 403.8818631413123
$ws.Range("C7").Value = 814.5043961409763
$ws.Range("E7").Value = 1732.2522243302978
$ws.Range("G7").Value = 874.550102548951
$ws.Range("J7").Value = 538.049241798567
$ws.Range("C8").Value = 902.8117247435255
$ws.Range("E8").Value = 1755.6532004432918
$ws.Range("G8").Value = 914.2163279336471
$ws.Range("J8").Value = 606.4236026044599
$ws.Range("C9").Value = 964.2669056735263
$ws.Range("E9").Value = 1823.0980967187024
$ws.Range("G9").Value = 989.2700660622756
$ws.Range("J9").Value = 1120.6318661902958
$ws.Range("C10").Value = 1081.3982042134373
$ws.Range("E10").Value = 1758.7518930502324
$ws.Range("G10").Value = 777.6856690849313
$ws.Range("J10").Value = 675.8701970175887
$ws.Range("C11").Value = 1196.5457985170783
$ws.Range("E11").Value = 1736.4472428954737
$ws.Range("G11").Value = 583.147641217855
$ws.Range("J11").Value = 537.9952057606333
